# Update cryptocurrency price/volume data per upstream refresh (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.381.47"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "3.947.04"
$ws.Range("E3").Value = "  +4.50%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'488.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +10.27%  "

$ws.Range("D6").Value = "'148.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.76%  "

$ws.Range("D7").Value = "'0.627"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.38%  "

$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'0.731"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("E10").Value = "  +11.57%  "

$ws.Range("E11").Value = "  +14.80%  "

$ws.Range("D12").Value = "'43.15"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "4.578.17"
$ws.Range("E13").Value = "  +4.80%  "

$ws.Range("D14").Value = "'10.41"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").Value = "'15.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.17%  "

$ws.Range("D16").Value = "3.974.72"
$ws.Range("E16").Value = "  +6.30%  "

$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").Value = "'19.98"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("D20").Value = "67.500.49"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("D21").Value = "'433.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.39%  "

$ws.Range("E22").Value = "  +4.71%  "

$ws.Range("D23").Value = "'14.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "'87.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.77%  "

$ws.Range("D25").Value = "'3.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.24%  "

$ws.Range("D26").Value = "'38.53"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.71%  "

$ws.Range("D27").Value = "'9.87"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").Value = "'10.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.29%  "

$ws.Range("D29").Value = "'722.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("D30").Value = "'0.132"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.27%  "

$ws.Range("D31").Value = "'13.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("D33").Value = "'42.33"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("D34").Value = "0.0₃0853"
$ws.Range("E34").Value = "  +27.99%  "

$ws.Range("D35").Value = "'58.73"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.20%  "

$ws.Range("D36").Value = "'0.154"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.99%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'5.42"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'0.0475"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("D40").Value = "'3.06"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.07%  "

$ws.Range("E41").Value = "  +3.40%  "

$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("E43").Value = "  +8.05%  "

$ws.Range("D44").Value = "'2.85"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.81%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("E47").Value = "  +4.56%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.11%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'148.42"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.85%  "

$ws.Range("D50").Value = "'2.90"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("D51").Value = "'25.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.59%  "
